# Release Documentation Input Variables.xlsx
# PR 2788: RelDoc update for PRD Release 2020.1 (PRD Release on Jan 23rd)

$wb = $excel.ActiveWorkbook

$wsQA   = $wb.Worksheets.Item("QA")
$wsProd = $wb.Worksheets.Item("Prod")

# ---------------------------------------------------------------------------
# Update the "Prod" sheet's input values for the new PRD 2020.1 release
# ---------------------------------------------------------------------------
$wsProd.Range("C2").Value  = "CHG0042994"                               # ChangeTicketNum
$wsProd.Range("C7").Value  = 43853                                      # DeployDate (2020-01-23)
$wsProd.Range("C8").Value  = 2020                                       # ReleaseYear
$wsProd.Range("C9").Value  = "1"                                        # ReleaseSeqNum
$wsProd.Range("C15").Value = "n"                                        # IconDbUpdated
$wsProd.Range("C22").Value = "RELEASE_2020_01"                          # IconReleaseBranchName
$wsProd.Range("C28").Value = "11.8"                                     # IrmaVer
$wsProd.Range("C29").Value = "11.8.0"                                   # IrmaLongVer
$wsProd.Range("C46").Value = "PublishTransferOrderService"              # TibcoAppsUpdated
$wsProd.Range("C47").Value = "FL,MA,MW,NA,NC,NE,PN,RM,SO,SP,SW"         # HcTibcoRegions
$wsProd.Range("C50").Value = "11.7.0"                                   # IrmaSuitePreviousVersion
$wsProd.Range("C54").Value = "<ol><li>POS Push Job<li>PeopleSoft Upload Job<li>PeopleSoft Transfer Upload Job<li>MILD SSIS Jobs</ol>"  # AdditionalIrmaComponentDeploySection
$wsProd.Range("C55").Value = "POS Push, PS Upload, PS Transfer Upload"  # AdditionalIrmaComponentList
$wsProd.Range("C57").Value = "<hr>IRMA Client<br>Icon API Controller<br>Mammoth DB<br>Mammoth Audit Service<br>Mammoth Hierarchy Class Listener<br>Mammoth Product Listener<br>Mammoth Item Locale Controller<br>Mammoth Web Support"  # AzureReleaseList

# Row 57 text grew (more <hr>/<br> segments) -- widen the row to fit, as in the target
$wsProd.Rows.Item(57).RowHeight = 90

# ---------------------------------------------------------------------------
# Flip the active tab from "QA" to "Prod" and restore each sheet's selection
# ---------------------------------------------------------------------------
$wsProd.Activate()
$wsProd.Range("B61").Select()
$excel.ActiveWindow.ScrollRow = 59

$wsQA.Activate()
$wsQA.Range("C38").Select()
$excel.ActiveWindow.ScrollRow = 2

$wsProd.Activate()
